# Insert a new weekly price record as row 334 in the "Fruta, Vega Modelo de
# Temuco - Mango" sheet. This shifts the existing rows 334-406 down to
# 335-407 (preserving all of their data and formatting), and fills the
# freshly inserted row 334 with the new record, copying the unchanging
# "template" columns from the row directly below (which now holds what used
# to be row 334) and then overwriting the columns that actually carry new
# data (Fecha, Volumen, Precio minimo/maximo/promedio, Origen, Precio $/Kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 334:406 down to 335:407, inserting a blank row 334.
$ws.Rows("334").Insert()

# Columns that stay the same as the template row immediately below (335)
# after the insert: A, B, C, E, F, G, H, I, J, K, L, Q, T.
$templateRow = 335
$newRow = 334
$sameCols = @("A","B","C","E","F","G","H","I","J","K","L","Q","T")
foreach ($col in $sameCols) {
    $ws.Range("$col$newRow").Value = $ws.Range("$col$templateRow").Value()
}

# (Inserting the row already carried the date-style formatting down from
# the row above onto D334, matching the other Fecha cells.)

# New record values.
$ws.Range("D$newRow").Value = 44798
$ws.Range("M$newRow").Value = 250
$ws.Range("N$newRow").Value = 10000
$ws.Range("O$newRow").Value = 10000
$ws.Range("P$newRow").Value = 10000
$ws.Range("R$newRow").Value = "Brasil"
$ws.Range("S$newRow").Value = 2500
